$d = $word.ActiveDocument

# 1) Remove the _GoBack bookmark currently sitting at the start of the
#    "Contact Amira and make plans" paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Put the new sentence into the trailing empty paragraph (the last
#    paragraph in the document, after "Figure out health insurance plan").
$paragraphs = $d.Paragraphs
$lastPara = $paragraphs.Item($paragraphs.Count)
$lastPara.Range.Text = "I EDITTED SOMETHING IN THIS FILE "

# 3) Re-create the _GoBack bookmark at the end of that paragraph's text
#    (after the run, before the paragraph mark). A temporary placeholder
#    character is used so the bookmark collapses cleanly after the run
#    instead of spanning the whole paragraph, then the placeholder is
#    removed.
$newLastPara = $paragraphs.Item($paragraphs.Count)
$endRange = $newLastPara.Range.Duplicate()
$endRange.Collapse(0)
$endRange.InsertAfter("X")
$endRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $endRange)

$finalPara = $paragraphs.Item($paragraphs.Count)
$placeholder = $finalPara.Range.Duplicate()
[void]$placeholder.MoveEnd(1, -1)
[void]$placeholder.MoveStart(1, $placeholder.End - $placeholder.Start - 1)
$placeholder.Text = ""
